$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.032.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.910.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3809"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07358"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9340"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07811"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.910.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.505"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.634"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.89%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008871"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "28.056.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.170"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.133.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.912"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.72%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.71%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.973"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08957"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.281"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.258"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.71%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7752"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.663"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.616"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5535"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05296"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.993"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.031"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aptos"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.518"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Algorand"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1529"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "108.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.22%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4834"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.651"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06073"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.21%  "
$ws.Range("E51").Style = "Normal"

